# Apply leave registry corrections for George Smith 2020 leave registry sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("George Smith_2020_leave_registry")

# Header date fields
$ws.Range("G8").Value = "June 03, 2020"
$ws.Range("G10").Value = "June 03, 2020"

# Vacation leave summary row (13)
$ws.Range("G13").Value = 5.625
$ws.Range("I13").Value = 5.125

# Sick leave summary row (14)
$ws.Range("G14").Value = 3.744
$ws.Range("H14").Value = 4.5
$ws.Range("I14").Value = -0.7560000000000002

# Details row 19 - updated date, day count, and leave type
$ws.Range("F19").Value = "May 19, 2020"
$ws.Range("G19").Value = 4.5
$ws.Range("I19").Value = "SL"

# Details row 20 - newly populated entry
$ws.Range("F20").Value = "May 20, 2020"
$ws.Range("G20").Value = 0.5
$ws.Range("I20").Value = "VL"

# Row 24 bi-monthly earned credits
$ws.Range("C24").Value = 0.625
$ws.Range("D24").Value = 0.416

# Totals row 40
$ws.Range("C40").Value = 5.625
$ws.Range("D40").Value = 3.744

# Offense log row 45
$ws.Range("F45").Value = "April 13, 2020"
$ws.Range("G45").Value = "No Time-in"
